# COREESG_holdings.xlsx update:
#   - Roll the "as of" date in the confidential disclosure text (A10) from
#     2021-03-26 to 2021-03-29.
#   - Refresh the Weight (D) / Percent Change (E) figures for rows 2-7.
#
# The worksheet is protected (legacy password hash "D382"); the target
# cells are all locked by default, so a direct .Value write is rejected by
# the protection guard. Rather than calling Worksheet.Unprotect()/Protect()
# (which would re-hash any supplied password and replace the original
# <sheetProtection password="D382" .../> with a brand-new modern hash),
# we momentarily unlock just the cell we are writing, write it, then copy
# the number format back from an untouched neighboring cell that already
# carries the correct original style/locked state. That restores the cell
# to its original look while the sheet itself is never unprotected, so
# <sheetProtection password="D382" .../> round-trips unchanged.
#
# NOTE: named parameter binding (-Name value) is unreliable for typed
# function params in this host, so helper functions below use plain,
# untyped params invoked positionally.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-ProtectedValue($Sheet, $App, $TargetAddress, $Value, $FormatSourceAddress) {
    $target = $Sheet.Range($TargetAddress)
    $target.Locked = $false
    $target.Value = $Value

    # Restore the original number format / locked state by pasting the
    # formatting (only) from an untouched cell that already has the
    # correct original style.
    $Sheet.Range($FormatSourceAddress).Copy()
    $target.PasteSpecial(-4122)  # xlPasteFormats
    $App.CutCopyMode = $false
}

# --- A10: confidential disclosure text, date 2021-03-26 -> 2021-03-29 ---
$newText = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-29 for illustrative purposes only and are subject to change."
Set-ProtectedValue $ws $excel "A10" $newText "A11"
# Re-fit the row after the multi-line text edit so no stray custom height sticks.
$ws.Rows.Item(10).EntireRow.AutoFit()

# --- Row 2 (NULG) ---
Set-ProtectedValue $ws $excel "D2" 0.2455707393709217 "D3"
Set-ProtectedValue $ws $excel "E2" -0.00226559776925761 "E3"

# --- Row 3 (NULV) ---
Set-ProtectedValue $ws $excel "D3" 0.4991281292852269 "D4"
Set-ProtectedValue $ws $excel "E3" -0.001096791883740011 "E4"

# --- Row 4 (NUMG) ---
Set-ProtectedValue $ws $excel "D4" 0.09813809595567664 "D5"
Set-ProtectedValue $ws $excel "E4" -0.008474576271186529 "E5"

# --- Row 5 (NUMV) ---
Set-ProtectedValue $ws $excel "D5" 0.09959355113859196 "D6"
Set-ProtectedValue $ws $excel "E5" -0.00897770055024627 "E6"

# --- Row 6 (NUSC) ---
Set-ProtectedValue $ws $excel "D6" 0.05756948424958289 "D2"
Set-ProtectedValue $ws $excel "E6" -0.02341757477393924 "E2"

# --- Row 7 (Total) ---
Set-ProtectedValue $ws $excel "E7" -0.004177741760714948 "E2"
